$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(410.0, 396.0, 413.0, 392.0, 390.0, 419.0, 397.0, 394.0, 446.0, 432.0, 391.0, 404.0, 407.0, 393.0, 405.0, 410.0, 396.0, 413.0, 392.0, 390.0, 446.0, 419.0, 432.0, 397.0, 391.0, 394.0, 404.0, 407.0, 393.0, 405.0)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
